$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Insert a new column before column J (10). This shifts the existing
# "have_seizure_lbl" column (J) to K and "eeg_lbl_path" column (K) to L,
# while preserving their values/styles/comments, and grows the used range
# from A1:K32 to A1:L32.
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).Insert()

# Header for the newly inserted column
$ws.Range("J1").Value = "FluoroVolume"

# FluoroVolume data values (row 11 intentionally has no value, matching source)
$ws.Range("J2").Value = 16.8
$ws.Range("J3").Value = 9.4
$ws.Range("J4").Value = 25.2
$ws.Range("J5").Value = 26.2
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 14.3
$ws.Range("J8").Value = 12.3
$ws.Range("J9").Value = 14.3
$ws.Range("J10").Value = 2.5
$ws.Range("J12").Value = 4.5999999999999996
$ws.Range("J13").Value = 6

# Match the "0.00" number format used by the neighbouring SkullFluoro_mm2 column
$ws.Range("J2:J10").NumberFormat = "0.00"
$ws.Range("J12:J13").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Column widths: column J now takes the old "eeg_lbl_path" width, column K
# gets a new width, column L keeps the old "eeg_lbl_path" width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 31.42
$ws.Columns.Item(11).ColumnWidth = 14.88

# ---------------------------------------------------------------------------
# Fix up the hyperlinks on the "eeg_lbl_path" column. Inserting a column
# shifts cell values/styles automatically, but leaves hyperlink anchors
# pointing at their old (now incorrect) addresses, so re-create them.
# ---------------------------------------------------------------------------
$ws.Range("K2").Hyperlinks.Delete()
$ws.Range("K4").Hyperlinks.Delete()
$ws.Range("K5").Hyperlinks.Delete()
$ws.Range("K6").Hyperlinks.Delete()

# Remember the untouched hyperlink cell style so it can be re-applied after
# Hyperlinks.Add (which otherwise stamps in a freshly duplicated style).
$linkStyle = $ws.Range("L4").Style()

$ws.Hyperlinks.Add($ws.Range("L2"), "file:///\\neurodata\Lab%20Neurophysiology%20root\EEG%20conversion\NatySST_TdTET339")
$ws.Hyperlinks.Add($ws.Range("L4"), "file:///\\neurodata\Lab%20Neurophysiology%20root\EEG%20conversion\NatymTORET283")
$ws.Hyperlinks.Add($ws.Range("L5"), "file:///\\neurodata\Lab%20Neurophysiology%20root\EEG%20Naty\mTOR%20MUT\Naty%20SST_TdT%20ET%20343")
$ws.Hyperlinks.Add($ws.Range("L6"), "file:///\\neurodata\Lab%20Neurophysiology%20root\EEG%20Naty\mTOR%20MUT\Naty%20SST_TdT%20ET%20413")

$ws.Range("L2").Style = $linkStyle
$ws.Range("L4").Style = $linkStyle
$ws.Range("L5").Style = $linkStyle
$ws.Range("L6").Style = $linkStyle

# ---------------------------------------------------------------------------
# View state (best effort)
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J17").Select()

Write-Host "Edit applied"
